$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2450031936168671
$ws.Range("B1").Value = 0.3070769608020782
$ws.Range("C1").Value = 0.4467190504074097
$ws.Range("D1").Value = 2.051134347915649
$ws.Range("E1").Value = 5.570147037506104
